# Add a new "Sede Reason" column (O) to the inquiring-requests report,
# mirroring the header style of the last existing column (N, "Domain Reason"),
# then extend the sheet's AutoFilter and the workbook's
# _xlnm._FilterDatabase defined name to cover the new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the autofilter is off before we touch the header row / range,
# so re-enabling it afterwards rebuilds it over the full A1:O1 range.
$ws.AutoFilterMode = $false

# New header cell, copying the formatting (gray fill) used by the other
# header cells from its neighbour.
$ws.Range("N1").Copy()
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("O1").Value = "Sede Reason"

# Size the new column like the other bestFit header columns.
$ws.Columns.Item(15).ColumnWidth = 16.2

# Re-apply the autofilter across the now-wider header range.
$null = $ws.Range("A1:O1").AutoFilter()

# Widen the _xlnm._FilterDatabase defined name to match.
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name() -eq "Data!_FilterDatabase") {
        $n.RefersTo = "=Data!`$A`$1:`$O`$1"
    }
}

Write-Output "Added Sede Reason column"
